$wb = $excel.ActiveWorkbook

# Sheet ALC, row 31
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 37037132
$ws.Range("I31").Value = 37037132
$ws.Range("K31").Value = 111111396
$ws.Range("M31").Value = -111111166

# Sheet ALC, row 55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 28
$ws.Range("I55").Value = 28
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 28
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 186
$ws.Range("N55").ClearContents()

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3533.68
$ws.Range("J138").Value = 3494.8167
$ws.Range("L138").Value = 10484.4501
$ws.Range("N138").Value = -20764.4501

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3983.6667
$ws.Range("I141").Value = 3981.75
$ws.Range("J141").Value = 3999
$ws.Range("K141").Value = 11945.25
$ws.Range("L141").Value = 11997
$ws.Range("M141").Value = -6765.25
$ws.Range("N141").Value = -22357

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1742.2632
$ws.Range("I2").Value = 791.9091
$ws.Range("K2").Value = 791.9091
$ws.Range("M2").Value = -678.9091

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5155.4
$ws.Range("J61").Value = 7998.875
$ws.Range("L61").Value = 7998.875
$ws.Range("N61").Value = -8422.875

# Sheet ARM, row 101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 43144.75
$ws.Range("J101").Value = 43144.75
$ws.Range("L101").Value = 43144.75
$ws.Range("N101").Value = -49634.75

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1742.2632
$ws.Range("I116").Value = 791.9091
$ws.Range("K116").Value = 791.9091
$ws.Range("M116").Value = 1502.0909

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2143.6428
$ws.Range("I132").Value = 1922.8462
$ws.Range("K132").Value = 5768.5386
$ws.Range("M132").Value = -3238.5386

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5155.4
$ws.Range("J136").Value = 7998.875
$ws.Range("L136").Value = 23996.625
$ws.Range("N136").Value = -29096.625

# Sheet ARM, row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 23999.2
$ws.Range("I139").Value = 23999.2
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 23999.2
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -18859.2
$ws.Range("N139").ClearContents()

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1742.2632
$ws.Range("I3").Value = 791.9091
$ws.Range("K3").Value = 791.9091
$ws.Range("M3").Value = -677.9091

# Sheet BSM, row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 39253
$ws.Range("J81").Value = 39253
$ws.Range("L81").Value = 39253
$ws.Range("N81").Value = -41375

# Sheet BSM, row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 39253
$ws.Range("J84").Value = 39253
$ws.Range("L84").Value = 117759
$ws.Range("N84").Value = -128367

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 566
$ws.Range("I134").Value = 488
$ws.Range("J134").Value = 800
$ws.Range("K134").Value = 1464
$ws.Range("L134").Value = 2400
$ws.Range("M134").Value = 1071
$ws.Range("N134").Value = -7470

# Sheet CRP, row 15
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1007.5
$ws.Range("I15").Value = 1006
$ws.Range("J15").Value = 1009
$ws.Range("K15").Value = 1006
$ws.Range("L15").Value = 1009
$ws.Range("M15").Value = -836
$ws.Range("N15").Value = -1349

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2658.7144
$ws.Range("I16").Value = 2658.7144
$ws.Range("K16").Value = 2658.7144
$ws.Range("M16").Value = -2371.7144

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 27494.4
$ws.Range("I22").Value = 2990
$ws.Range("K22").Value = 2990
$ws.Range("M22").Value = -2640

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6421.625
$ws.Range("I86").Value = 7644.5
$ws.Range("J86").Value = 5198.75
$ws.Range("K86").Value = 7644.5
$ws.Range("L86").Value = 5198.75
$ws.Range("M86").Value = -6521.5
$ws.Range("N86").Value = -7444.75

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 6421.625
$ws.Range("I89").Value = 7644.5
$ws.Range("J89").Value = 5198.75
$ws.Range("K89").Value = 38222.5
$ws.Range("L89").Value = 25993.75
$ws.Range("M89").Value = -32606.5
$ws.Range("N89").Value = -37225.75

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2658.7144
$ws.Range("I113").Value = 2658.7144
$ws.Range("K113").Value = 2658.7144
$ws.Range("M113").Value = -488.7143999999998

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1643.5454
$ws.Range("I132").Value = 1355.25
$ws.Range("K132").Value = 4065.75
$ws.Range("M132").Value = -1535.75

# Sheet CUL, row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Sheet CUL, row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Sheet CUL, row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3560.6365
$ws.Range("J137").Value = 4134.75
$ws.Range("L137").Value = 12404.25
$ws.Range("N137").Value = -22604.25

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1769.5
$ws.Range("I132").Value = 1769.5
$ws.Range("K132").Value = 5308.5
$ws.Range("M132").Value = -2778.5

# Sheet LTW, row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 360.54544
$ws.Range("I16").Value = 361.6
$ws.Range("K16").Value = 361.6
$ws.Range("M16").Value = -191.6

# Sheet LTW, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1158.7222
$ws.Range("I22").Value = 972.75
$ws.Range("J22").Value = 1307.5
$ws.Range("K22").Value = 972.75
$ws.Range("L22").Value = 1307.5
$ws.Range("M22").Value = -677.75
$ws.Range("N22").Value = -1897.5

# Sheet LTW, row 25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Sheet LTW, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1158.7222
$ws.Range("I27").Value = 972.75
$ws.Range("J27").Value = 1307.5
$ws.Range("K27").Value = 972.75
$ws.Range("L27").Value = 1307.5
$ws.Range("M27").Value = -865.75
$ws.Range("N27").Value = -1521.5

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3965.4614
$ws.Range("I40").Value = 2958.6667
$ws.Range("K40").Value = 2958.6667
$ws.Range("M40").Value = -2822.6667

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2346.0588
$ws.Range("I46").Value = 883.1429000000001
$ws.Range("J46").Value = 3370.1
$ws.Range("K46").Value = 883.1429000000001
$ws.Range("L46").Value = 3370.1
$ws.Range("M46").Value = -695.1429000000001
$ws.Range("N46").Value = -3746.1

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4015.8518
$ws.Range("I61").Value = 3712.238
$ws.Range("K61").Value = 3712.238
$ws.Range("M61").Value = -3510.238

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2200
$ws.Range("I68").Value = 600
$ws.Range("K68").Value = 600
$ws.Range("M68").Value = 149

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2200
$ws.Range("I71").Value = 600
$ws.Range("K71").Value = 3000
$ws.Range("M71").Value = 744

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4015.8518
$ws.Range("I113").Value = 3712.238
$ws.Range("K113").Value = 3712.238
$ws.Range("M113").Value = -1542.238

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5919.35
$ws.Range("I122").Value = 5786.931
$ws.Range("K122").Value = 17360.793
$ws.Range("M122").Value = -14910.793

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4013.1667
$ws.Range("I132").Value = 1914.8
$ws.Range("J132").Value = 14505
$ws.Range("K132").Value = 5744.4
$ws.Range("L132").Value = 43515
$ws.Range("M132").Value = -3214.4
$ws.Range("N132").Value = -48575

# Sheet WVR, row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 24376
$ws.Range("I45").Value = 20635.666
$ws.Range("J45").Value = 29986.5
$ws.Range("K45").Value = 20635.666
$ws.Range("L45").Value = 29986.5
$ws.Range("M45").Value = -20144.666
$ws.Range("N45").Value = -30968.5

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2556.9167
$ws.Range("I132").Value = 2268.8
$ws.Range("J132").Value = 3997.5
$ws.Range("K132").Value = 6806.400000000001
$ws.Range("L132").Value = 11992.5
$ws.Range("M132").Value = -4276.400000000001
$ws.Range("N132").Value = -17052.5

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4960.5
$ws.Range("I136").Value = 5045.1113
$ws.Range("K136").Value = 15135.3339
$ws.Range("M136").Value = -12585.3339
